$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 103.4275383333333
$ws.Cells.Item(2, 8).Value = 310.282615
$ws.Cells.Item(2, 9).Value = 0.2485530285127421
$ws.Cells.Item(2, 10).Value = 0.2485530285127421
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 8.131233999999999
$ws.Cells.Item(2, 14).Value = 24.393702
$ws.Cells.Item(2, 15).Value = 0.02090995573015822
$ws.Cells.Item(2, 16).Value = 0.02090995573015823
$ws.Cells.Item(2, 17).Value = 840.9935162323034
$ws.Cells.Item(2, 18).Value = 7568.94164609073
$ws.Cells.Item(2, 19).Value = 0.005197232822798192
$ws.Cells.Item(2, 20).Value = 0.005197232822798192

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 103.4275383333333
$ws.Cells.Item(3, 8).Value = 310.282615
$ws.Cells.Item(3, 9).Value = 0.2485530285127421
$ws.Cells.Item(3, 10).Value = 0.2485530285127421
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 243.3763986666667
$ws.Cells.Item(3, 14).Value = 730.1291960000001
$ws.Cells.Item(3, 15).Value = 0.625857000534647
$ws.Cells.Item(3, 16).Value = 0.6258570005346471
$ws.Cells.Item(3, 17).Value = 25171.82180252529
$ws.Cells.Item(3, 18).Value = 226546.3962227276
$ws.Cells.Item(3, 19).Value = 0.1555586528987874
$ws.Cells.Item(3, 20).Value = 0.1555586528987874

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 103.4275383333333
$ws.Cells.Item(4, 8).Value = 310.282615
$ws.Cells.Item(4, 9).Value = 0.2485530285127421
$ws.Cells.Item(4, 10).Value = 0.2485530285127421
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 103.9426383333333
$ws.Cells.Item(4, 14).Value = 311.827915
$ws.Cells.Item(4, 15).Value = 0.2672947262403034
$ws.Cells.Item(4, 16).Value = 0.2672947262403035
$ws.Cells.Item(4, 17).Value = 10750.53121068864
$ws.Cells.Item(4, 18).Value = 96754.78089619774
$ws.Cells.Item(4, 19).Value = 0.06643691371251173
$ws.Cells.Item(4, 20).Value = 0.06643691371251174

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 103.4275383333333
$ws.Cells.Item(5, 8).Value = 310.282615
$ws.Cells.Item(5, 9).Value = 0.2485530285127421
$ws.Cells.Item(5, 10).Value = 0.2485530285127421
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 33.41874933333333
$ws.Cells.Item(5, 14).Value = 100.256248
$ws.Cells.Item(5, 15).Value = 0.08593831749489127
$ws.Cells.Item(5, 16).Value = 0.08593831749489128
$ws.Cells.Item(5, 17).Value = 3456.418977725391
$ws.Cells.Item(5, 18).Value = 31107.77079952852
$ws.Cells.Item(5, 19).Value = 0.02136022907864479
$ws.Cells.Item(5, 20).Value = 0.02136022907864479

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 216.130539
$ws.Cells.Item(6, 8).Value = 648.391617
$ws.Cells.Item(6, 9).Value = 0.5193964865470273
$ws.Cells.Item(6, 10).Value = 0.5193964865470272
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 8.131233999999999
$ws.Cells.Item(6, 14).Value = 24.393702
$ws.Cells.Item(6, 15).Value = 0.02090995573015822
$ws.Cells.Item(6, 16).Value = 0.02090995573015823
$ws.Cells.Item(6, 17).Value = 1757.407987155126
$ws.Cells.Item(6, 18).Value = 15816.67188439613
$ws.Cells.Item(6, 19).Value = 0.01086055754009806
$ws.Cells.Item(6, 20).Value = 0.01086055754009806

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 216.130539
$ws.Cells.Item(7, 8).Value = 648.391617
$ws.Cells.Item(7, 9).Value = 0.5193964865470273
$ws.Cells.Item(7, 10).Value = 0.5193964865470272
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 243.3763986666667
$ws.Cells.Item(7, 14).Value = 730.1291960000001
$ws.Cells.Item(7, 15).Value = 0.625857000534647
$ws.Cells.Item(7, 16).Value = 0.6258570005346471
$ws.Cells.Item(7, 17).Value = 52601.07222370555
$ws.Cells.Item(7, 18).Value = 473409.65001335
$ws.Cells.Item(7, 19).Value = 0.3250679271585566
$ws.Cells.Item(7, 20).Value = 0.3250679271585566

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 216.130539
$ws.Cells.Item(8, 8).Value = 648.391617
$ws.Cells.Item(8, 9).Value = 0.5193964865470273
$ws.Cells.Item(8, 10).Value = 0.5193964865470272
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 103.9426383333333
$ws.Cells.Item(8, 14).Value = 311.827915
$ws.Cells.Item(8, 15).Value = 0.2672947262403034
$ws.Cells.Item(8, 16).Value = 0.2672947262403035
$ws.Cells.Item(8, 17).Value = 22465.1784480654
$ws.Cells.Item(8, 18).Value = 202186.6060325886
$ws.Cells.Item(8, 19).Value = 0.1388319416817631
$ws.Cells.Item(8, 20).Value = 0.1388319416817631

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 216.130539
$ws.Cells.Item(9, 8).Value = 648.391617
$ws.Cells.Item(9, 9).Value = 0.5193964865470273
$ws.Cells.Item(9, 10).Value = 0.5193964865470272
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 33.41874933333333
$ws.Cells.Item(9, 14).Value = 100.256248
$ws.Cells.Item(9, 15).Value = 0.08593831749489127
$ws.Cells.Item(9, 16).Value = 0.08593831749489128
$ws.Cells.Item(9, 17).Value = 7222.812306119224
$ws.Cells.Item(9, 18).Value = 65005.31075507301
$ws.Cells.Item(9, 19).Value = 0.04463606016660945
$ws.Cells.Item(9, 20).Value = 0.04463606016660945

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 71.607325
$ws.Cells.Item(10, 8).Value = 214.821975
$ws.Cells.Item(10, 9).Value = 0.1720839321833696
$ws.Cells.Item(10, 10).Value = 0.1720839321833696
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 8.131233999999999
$ws.Cells.Item(10, 14).Value = 24.393702
$ws.Cells.Item(10, 15).Value = 0.02090995573015822
$ws.Cells.Item(10, 16).Value = 0.02090995573015823
$ws.Cells.Item(10, 17).Value = 582.25591568905
$ws.Cells.Item(10, 18).Value = 5240.30324120145
$ws.Cells.Item(10, 19).Value = 0.003598267403825807
$ws.Cells.Item(10, 20).Value = 0.003598267403825808

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 71.607325
$ws.Cells.Item(11, 8).Value = 214.821975
$ws.Cells.Item(11, 9).Value = 0.1720839321833696
$ws.Cells.Item(11, 10).Value = 0.1720839321833696
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 243.3763986666667
$ws.Cells.Item(11, 14).Value = 730.1291960000001
$ws.Cells.Item(11, 15).Value = 0.625857000534647
$ws.Cells.Item(11, 16).Value = 0.6258570005346471
$ws.Cells.Item(11, 17).Value = 17427.53287665357
$ws.Cells.Item(11, 18).Value = 156847.7958898821
$ws.Cells.Item(11, 19).Value = 0.1076999336364913
$ws.Cells.Item(11, 20).Value = 0.1076999336364913

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 71.607325
$ws.Cells.Item(12, 8).Value = 214.821975
$ws.Cells.Item(12, 9).Value = 0.1720839321833696
$ws.Cells.Item(12, 10).Value = 0.1720839321833696
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 103.9426383333333
$ws.Cells.Item(12, 14).Value = 311.827915
$ws.Cells.Item(12, 15).Value = 0.2672947262403034
$ws.Cells.Item(12, 16).Value = 0.2672947262403035
$ws.Cells.Item(12, 17).Value = 7443.054284492458
$ws.Cells.Item(12, 18).Value = 66987.48856043213
$ws.Cells.Item(12, 19).Value = 0.04599712754330871
$ws.Cells.Item(12, 20).Value = 0.04599712754330872

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 71.607325
$ws.Cells.Item(13, 8).Value = 214.821975
$ws.Cells.Item(13, 9).Value = 0.1720839321833696
$ws.Cells.Item(13, 10).Value = 0.1720839321833696
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 33.41874933333333
$ws.Cells.Item(13, 14).Value = 100.256248
$ws.Cells.Item(13, 15).Value = 0.08593831749489127
$ws.Cells.Item(13, 16).Value = 0.08593831749489128
$ws.Cells.Item(13, 17).Value = 2393.027244605533
$ws.Cells.Item(13, 18).Value = 21537.2452014498
$ws.Cells.Item(13, 19).Value = 0.01478860359974375
$ws.Cells.Item(13, 20).Value = 0.01478860359974375

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 24.953198
$ws.Cells.Item(14, 8).Value = 74.859594
$ws.Cells.Item(14, 9).Value = 0.05996655275686102
$ws.Cells.Item(14, 10).Value = 0.05996655275686102
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 8.131233999999999
$ws.Cells.Item(14, 14).Value = 24.393702
$ws.Cells.Item(14, 15).Value = 0.02090995573015822
$ws.Cells.Item(14, 16).Value = 0.02090995573015823
$ws.Cells.Item(14, 17).Value = 202.900291986332
$ws.Cells.Item(14, 18).Value = 1826.102627876988
$ws.Cells.Item(14, 19).Value = 0.001253897963436161
$ws.Cells.Item(14, 20).Value = 0.001253897963436162

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 24.953198
$ws.Cells.Item(15, 8).Value = 74.859594
$ws.Cells.Item(15, 9).Value = 0.05996655275686102
$ws.Cells.Item(15, 10).Value = 0.05996655275686102
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 243.3763986666667
$ws.Cells.Item(15, 14).Value = 730.1291960000001
$ws.Cells.Item(15, 15).Value = 0.625857000534647
$ws.Cells.Item(15, 16).Value = 0.6258570005346471
$ws.Cells.Item(15, 17).Value = 6073.01946445627
$ws.Cells.Item(15, 18).Value = 54657.17518010643
$ws.Cells.Item(15, 19).Value = 0.03753048684081171
$ws.Cells.Item(15, 20).Value = 0.03753048684081171

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 24.953198
$ws.Cells.Item(16, 8).Value = 74.859594
$ws.Cells.Item(16, 9).Value = 0.05996655275686102
$ws.Cells.Item(16, 10).Value = 0.05996655275686102
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 103.9426383333333
$ws.Cells.Item(16, 14).Value = 311.827915
$ws.Cells.Item(16, 15).Value = 0.2672947262403034
$ws.Cells.Item(16, 16).Value = 0.2672947262403035
$ws.Cells.Item(16, 17).Value = 2593.701234974057
$ws.Cells.Item(16, 18).Value = 23343.31111476651
$ws.Cells.Item(16, 19).Value = 0.01602874330271988
$ws.Cells.Item(16, 20).Value = 0.01602874330271988

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 24.953198
$ws.Cells.Item(17, 8).Value = 74.859594
$ws.Cells.Item(17, 9).Value = 0.05996655275686102
$ws.Cells.Item(17, 10).Value = 0.05996655275686102
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 33.41874933333333
$ws.Cells.Item(17, 14).Value = 100.256248
$ws.Cells.Item(17, 15).Value = 0.08593831749489127
$ws.Cells.Item(17, 16).Value = 0.08593831749489128
$ws.Cells.Item(17, 17).Value = 833.9046690270346
$ws.Cells.Item(17, 18).Value = 7505.142021243312
$ws.Cells.Item(17, 19).Value = 0.00515342464989327
$ws.Cells.Item(17, 20).Value = 0.00515342464989327
